# Applies the documented edits to docs/安装测试文档.docx
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) After the "Nginx" hyperlink bullet, append a plain-text run with
#    install instructions.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Nginx") {
        $insertPoint = $p.Range.End - 1
        $r = $d.Range($insertPoint, $insertPoint)
        $r.InsertAfter("（下载Stable Version的nginx/Windows，解压压缩包即可）")
    }
}

# ---------------------------------------------------------------------
# 2) "双击ManagerBackend.exe运行项目" -> "双击ManagerBackend.exe运行该子项目"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "双击ManagerBackend.exe运行项目，记录下监听IP", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "双击ManagerBackend.exe运行该子项目，记录下监听IP", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "中的目录，替换为前端文件解压后的目录" gains a trailing 。, a space and
#    a new bold+italic warning sentence (split across two runs).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "中的目录，替换为前端文件解压后的目录") {
        $full = "中的目录，替换为前端文件解压后的目录。 请将Windows默认的反斜杠/，否则Nginx运行报错"
        $p.Range.Text = $full
        $start = $p.Range.Start
        $prefixLen = ("中的目录，替换为前端文件解压后的目录。 ").Length
        $midLen = ("请将Windows默认的反斜杠").Length

        $rb1 = $d.Range($start + $prefixLen, $start + $prefixLen + $midLen)
        $rb1.Font.Bold = 1
        $rb1.Font.Italic = 1

        $rb2 = $d.Range($start + $prefixLen + $midLen, $start + $full.Length)
        $rb2.Font.Bold = 1
        $rb2.Font.Italic = 1
    }
}

# ---------------------------------------------------------------------
# 4) "（注意目录、网址后的分号，丢失会导致Nginx启动失败）" -> reworded and
#    the whole sentence becomes bold+italic.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "（注意目录、网址后的分号，丢失会导致Nginx启动失败）") {
        $p.Range.Text = "注意 目录、网址 后的分号，丢失会导致Nginx启动失败"
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Font.Bold = 1
        $r.Font.Italic = 1
    }
}

# ---------------------------------------------------------------------
# 5) "在Ngnix根目录打开终端，输入" gains parenthetical PowerShell hint.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "在Ngnix根目录打开终端，输入", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "在Ngnix根目录打开终端（在目录的空白处按住Shift键，右键，选择“在此处打开PowerShell窗口”），输入",
    2) | Out-Null

# ---------------------------------------------------------------------
# 6) "启动项目" gains extra trailing guidance.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "启动项目", $true, $false, $false, $false, $false, $true, 1, $false,
    "启动项目，Nginx窗口一闪而过说明启动成功，即可关闭终端窗口", 2) | Out-Null

# ---------------------------------------------------------------------
# 7) "至此项目开始运行，监听域名为第14、15行" -> "...观察配置文件第14、15行"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "至此项目开始运行，监听域名为第14、15行", $true, $false, $false, $false,
    $false, $true, 1, $false, "至此项目开始运行，观察配置文件第14、15行", 2) | Out-Null

# ---------------------------------------------------------------------
# 8) "，即http://localhost:8888/" gains trailing sentence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "，即http://localhost:8888/", $true, $false, $false, $false, $false,
    $true, 1, $false, "，即http://localhost:8888/，使用此网址即可访问项目", 2) | Out-Null
